$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3
$ws.Range("B2").Value = 15.6
$ws.Range("C2").Value = 83.3

$ws.Range("B3").Value = 15.8
$ws.Range("C3").Value = 83.2
$ws.Range("D3").Value = 26.3

# Add new rows 4 and 5
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = 13.8
$ws.Range("C4").Value = 83.2
$ws.Range("D4").Value = 26.3

$ws.Range("A5").Value = 18
$ws.Range("B5").Value = 15.3
$ws.Range("C5").Value = 83.6
$ws.Range("D5").Value = 26.3
